# Swap the contents of columns B:AD between two pairs of rows.
# Column A (sequential id) stays untouched in each row.
# Pair 1: rows 114 and 115
# Pair 2: rows 162 and 163

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$col = 2
while ($col -le 30) {
    $tmp = $ws.Cells.Item(114, $col).Value2
    $ws.Cells.Item(114, $col).Value = $ws.Cells.Item(115, $col).Value2
    $ws.Cells.Item(115, $col).Value = $tmp
    $col = $col + 1
}

$col = 2
while ($col -le 30) {
    $tmp = $ws.Cells.Item(162, $col).Value2
    $ws.Cells.Item(162, $col).Value = $ws.Cells.Item(163, $col).Value2
    $ws.Cells.Item(163, $col).Value = $tmp
    $col = $col + 1
}
